# Update the "Förändrad" column (C) for all data rows (2-18) from
# 2023-10-22 (45221) to 2023-10-25 (45224), reflecting the automatic
# update of the logging overview for VINDELN.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 18 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
